$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.78"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'22.41"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'5.474"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.05617"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'6.461"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'0.8047"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'1.041"
$ws.Range("D8").Style = "Normal"

$ws.Range("B9").Value = "'One"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.01164"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'8OneONEBestin24h"
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = "'WazirX"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1423"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'9WazirXWRX"
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.07303"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'10MandalaExchangeTokenMDX"
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.03175"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.02929"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'12BitrueCoinBTR"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'13BitMartTokenBMX"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.001661"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'14BitForexTokenBF"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "'MCDex"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'3.228"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'15MCDexMCB"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "'CoinExToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.04742"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'16CoinExTokenCET"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.006448"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'0.005066"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'0.001056"
$ws.Range("D20").Style = "Normal"

$ws.Range("D22").Value = "'3.986"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'3.381"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "'2.085"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'0.3319"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "'25ProBitTokenPROB"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0003307"
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").Value = "'0.04171"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006895"
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = "'CEJI"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.003508"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'41CEJICEJI"
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'BKEXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.1037"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'42BKEXTokenBKK"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.008530"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005658"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'0.6815"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'0.01542"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002105"
$ws.Range("D49").Style = "Normal"
